# Natmi following Dr Hou advice
# The Comp-Itgb3 edge table gains the "ECs" (endothelial cells) cluster as
# both a sending and a target cluster, turning the old 2-cluster x 2-cluster
# (minus self-pairs-with-itself-only-partially) layout into the full
# 3-cluster x 3-cluster grid (ECs / FAPs / sCs), i.e. rows 2-7 -> rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous data rows (2:7); header row 1 is left untouched.
$ws.Range("A2:T7").Clear()

# Sending cluster (col A) / Target cluster (col D) for each of the 9 rows
# of the 3x3 cluster grid. Ligand symbol (col B) is always "Comp" and
# Receptor symbol (col C) is always "Itgb3" for this Comp-Itgb3 sheet.
$clusters = @("ECs", "FAPs", "sCs")

$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ws.Range("A$row").Value = $sending
        $ws.Range("B$row").Value = "Comp"
        $ws.Range("C$row").Value = "Itgb3"
        $ws.Range("D$row").Value = $target
        $row++
    }
}

# Numeric edge-statistics columns (E:T) for rows 2-10, recomputed for the
# new 3x3 cluster grid.
$values = @{
    2  = @(2, 0.6666666666666666, 0.569336, 1.708008, 0.01936485544401529, 0.01936485544401529, 3, 1, 7.214110666666667, 21.642332, 0.4688823795981188, 0.4688823795981188, 4.107252910517333, 36.965276194656, 0.009079839501163478, 0.009079839501163478)
    3  = @(2, 0.6666666666666666, 0.569336, 1.708008, 0.01936485544401529, 0.01936485544401529, 3, 1, 7.110350666666666, 21.331052, 0.4621384803214003, 0.4621384803214003, 4.048178607157332, 36.433607464416, 0.008949244866540824, 0.008949244866540824)
    4  = @(2, 0.6666666666666666, 0.569336, 1.708008, 0.01936485544401529, 0.01936485544401529, 3, 1, 1.061296333333333, 3.183889, 0.06897914008048092, 0.06897914008048092, 0.6042342092346665, 5.438107883111999, 0.001335771076310995, 0.001335771076310995)
    5  = @(3, 1, 28.31465866666666, 84.94397599999999, 0.9630679809930072, 0.9630679809930072, 3, 1, 7.214110666666667, 21.642332, 0.4688823795981188, 0.4688823795981188, 204.2650811102257, 1838.385729992032, 0.4515656066427571, 0.4515656066427571)
    6  = @(3, 1, 28.31465866666666, 84.94397599999999, 0.9630679809930072, 0.9630679809930072, 3, 1, 7.110350666666666, 21.331052, 0.4621384803214003, 0.4621384803214003, 201.3271521269724, 1811.944369142752, 0.4450707731823076, 0.4450707731823076)
    7  = @(3, 1, 28.31465866666666, 84.94397599999999, 0.9630679809930072, 0.9630679809930072, 3, 1, 1.061296333333333, 3.183889, 0.06897914008048092, 0.06897914008048092, 30.05024342251821, 270.4521908026639, 0.06643160116794258, 0.06643160116794258)
    8  = @(3, 1, 0.516483, 1.549449, 0.01756716356297749, 0.01756716356297749, 3, 1, 7.214110666666667, 21.642332, 0.4688823795981188, 0.4688823795981188, 3.725965519452, 33.533689675068, 0.008236933454198252, 0.008236933454198252)
    9  = @(3, 1, 0.516483, 1.549449, 0.01756716356297749, 0.01756716356297749, 3, 1, 7.110350666666666, 21.331052, 0.4621384803214003, 0.4621384803214003, 3.672375243372, 33.051377190348, 0.008118462272551892, 0.008118462272551892)
    10 = @(3, 1, 0.516483, 1.549449, 0.01756716356297749, 0.01756716356297749, 3, 1, 1.061296333333333, 3.183889, 0.06897914008048092, 0.06897914008048092, 0.5481415141289999, 4.933273627161, 0.001211767836227344, 0.001211767836227344)
}

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in 2..10) {
    $data = $values[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $data[$i]
    }
}
